$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = '新しく お店を出さないでというテーマで話を始めていきます じゃあ 近い やる方いないですが 今回 ちょっと 僕するで何だか何だかいないですか 司会'
$ws.Range("C3").Value = '新しく お店を出すんなら何がいいか 僕まで 意見 出てないって意見 ある人としてください 意見てからでいい飲食店'
$ws.Range("C4").Value = '仕事で家から僕は カードショップではカードがそれに何で出したいと思いました 誰かいませんかね'
$ws.Range("C5").Value = '八木書店 英語である店を作るなら 高級料亭かな理由は 高級料亭 だから自分の店で接待できるしたら接待できるほど走行 給料って だからさ 位の高い人物とか誘うわけやん 誘ったお客でくるわけ やんまー こねって言うかま 話聞かせてもらえるかもしれんし'
$ws.Range("C6").Value = 'させてもらえるかもしれんしとなんか自分で高級持ってる 高級料亭 お好きな子 誘って 食いに行けるとかめちゃくちゃ 良くない 何か美味しいもん食えるし その なんか大衆向けな店ってほんまに 客室 選べん市 しょぼい タンクが低い店ってなんかほんまに 変なやつ多いからこうきそうだから 高給料って言っても金持ってるやつしか金からほんまに変な奴はこうへんわ まともなやつぐらいしか来なさそうやからかにお店が接客するしたとしてもそっちの方が楽 や し まくりやからそんな人数も一気にするわけでもないから 拡散も捌きやすいからお給料って'
$ws.Range("C7").Value = '白テニ 開いて高級なところにこだわったかな 今 お客の人で言えばカードショップはね もう やべえ奴らと来るオタクとか入ってないやつとか臭い もう使えば確実に属した高級料亭のがうまく再生しない 草なるしこのショップ やけど なんか 儲かりそうな 詐欺とか流行ってるし もう一つの理由の一つででしたがクソなんか料理とか奈良市 高の原 くって言っとるし なんか臭いけどそれがあり コレクター ちょっと退職した数'
$ws.Range("C8").Value = 'お薬 しますは人からも愛し合うようになりますか'
$ws.Range("C9").Value = '薬がマイナスなら今 今'
$ws.Range("C10").Value = '理由としてはなんかあんまりお酒の専門店あんまりいないようなイメージがある島 基本的なスーパーのコーナーと一緒に何か出しているところが多いからそういうお酒の専門店もないからそういうのだしたらいいかなと思いました 鍵の専門店である酒屋はその専門的ではないという どちらかと言うとは 地域によってはその近くに 専門家がなかったりとかするところもあるから まで来たらまあ その あんまりその辺 持ってか集まってない地域にお店だすっていう設定で一応 前もっててって思っておりあるよ'
$ws.Range("C11").Value = '鮭とかってほんまに地方とかに バラけてるよ 水がいいとか3時 近かったりとかそんなスーパーマーケットとかなんとかもお酒の種類 めちゃくちゃ 漆山 持ってきちゃったらやっぱり 酒屋 地酒 高尾 ゲット ろ剣とか それもう一度 県に一個か 2個か 3個ぐらい 絶対 酒屋あるからちょっとなんか 例えばあったらこの店は ワイン 専門とか 日本酒専門 とか なんか 種類が限られるようなことが多いと思うから だから何とか'
$ws.Range("C12").Value = '3月とかいろんなやつなんか 比べる 神できるような感じにすれば'
$ws.Range("C13").Value = '今飯 まとめるまとめ とかあんの これ これはなんかさっき出た新車 やと 俺のコレクターカードカードや 野菜だけのビーガンと高級料理とお酒の料理店で待ってないってのは思うの 多分好きなものやからみんな あれやろしまさき君が言ってた これって一番を止めた 終わったも何かあるんなら そのポイントとしてその水だから昨日にかけておいたとしたら'
$ws.Range("C14").Value = '今問題 新車としたらもう6月9日 そんなもん6-6 なきゃ こないよね この辺で何か美味しいお店は多分 ビーガン ビーガン ビーガンと俺のカードは高級料理店へ'
$ws.Range("C15").Value = 'チョッパーの客室っていうポイント ならまあ 働きでは思うなあ れあるやろ 市 高級料理店がちょっとあのそういう そういった点とか情報って書いたりできるしていうのもありで一番気持ちいいのか 狂うか なんか一番そのりって感じていいですか もう30分やべやべ 今からちょっと考え'
$ws.Range("C16").Value = '天 専門店じゃないから どちらかというとは誰やったらお酒飲む人が多いと思うので やはり今近くに 専門店とか地元に えっと どちらかと言うと住んでるところから近くに 専門店があったら通いやすいと思うので まあ どちらかというとお酒の専門店が身近にあったらいいかなーって第三者として第3と第三者と変なフォローを入れたところでどこの あれは 会議は'

$ws.Range("C17").Value = ""
$ws.Rows("18:20").Delete()
